$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the trailing comma after "Abrupt" ("Abrupt," -> "Abrupt")
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Abrupt,", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Abrupt", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from the last (empty) paragraph
#    to the end of the "Claudiquer" paragraph, right after the final
#    character of its text (matching the target OOXML).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}

# Locate the paragraph whose text is "Claudiquer" (Range.Text includes the
# trailing paragraph mark, so trim it before comparing).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "Claudiquer") {
        $targetPara = $p
        break
    }
}

if ($targetPara -ne $null) {
    $pr = $targetPara.Range
    # One-character range covering the last character of the paragraph's
    # text (just before the paragraph mark). Re-inserting this character
    # together with the bookmark markers plants the bookmark right after
    # it, without introducing a spurious extra paragraph break.
    $lastCharRange = $d.Range($pr.End - 2, $pr.End - 1)
    $lastChar = $lastCharRange.Text

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p>' + `
        '<w:r><w:rPr>' + `
        '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
        '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="fr-FR"/>' + `
        '</w:rPr><w:t xml:space="preserve">' + $lastChar + '</w:t></w:r>' + `
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
        '</w:p></w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $lastCharRange.InsertXML($xml)
}
